$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-6, columns D, L, M, N, O, P, Q, R, S, T
# (derived from a cyclic reshuffle of the original rows' data)
$data = @{
    2 = @{ D = 44344; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 750;  T = 18 }
    3 = @{ D = 44334; L = "Primera"; M = 120; N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada"; R = "Región de O'Higgins"; S = 1042; T = 12 }
    4 = @{ D = 44316; L = "Primera"; M = 60;  N = 17500; O = 18000; P = 17750; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins"; S = 1109; T = 16 }
    5 = @{ D = 44316; L = "Segunda"; M = 40;  N = 16000; O = 16000; P = 16000; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins"; S = 1000; T = 16 }
    6 = @{ D = 44330; L = "Primera"; M = 60;  N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 861;  T = 18 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 4).Value  = $row.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $row.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $row.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $row.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $row.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $row.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $row.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $row.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row.T   # T: Kg / unidad
}
